$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.472.26"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.912.80"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'245.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("D6").Value = "'0.9989"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.4815"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("D8").Value = "'0.2891"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.29%  "
$ws.Range("D9").Value = "'0.06722"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").Value = "'110.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.03%  "
$ws.Range("D11").Value = "'19.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.95%  "
$ws.Range("D12").Value = "1.909.42"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").Value = "'5.266"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").Value = "'0.6721"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").Value = "'287.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Value = "30.485.11"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000007599"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "'0.9989"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "'12.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("D21").Value = "2.163.47"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "'5.471"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.86%  "
$ws.Range("D23").Value = "'0.9986"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").Value = "'6.411"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.42%  "
$ws.Range("D25").Value = "'9.463"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("D26").Value = "'163.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("D27").Value = "'20.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.70%  "
$ws.Range("D28").Value = "'2.113"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("D29").Value = "'0.1054"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.23%  "
$ws.Range("E30").Value = "  +2.61%  "
$ws.Range("D31").Value = "'4.173"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").Value = "'4.044"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.47%  "
$ws.Range("D33").Value = "'0.04978"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("D34").Value = "'0.7284"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.09%  "
$ws.Range("D35").Value = "'1.133"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("D36").Value = "'0.9989"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").Value = "'0.02033"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.22%  "
$ws.Range("D39").Value = "'2.667"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").Value = "'110.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("D41").Value = "'2.012"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.20%  "
$ws.Range("D42").Value = "'0.4431"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.75%  "
$ws.Range("D43").Value = "'0.8674"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").Value = "'5.820"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("D45").Value = "'0.9990"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "'68.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").Value = "'7.334"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("D48").Value = "'48.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.02%  "
$ws.Range("D49").Value = "'9.307"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").Value = "'0.1240"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("D51").Value = "'34.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.45%  "
